# edit.ps1 - apply the appendix-z.docx changes described in the diff:
#  1. Bump the TDD Cookbook Development version from 1.0.0 to 1.1.0, and split the
#     run so "CentOS" (a word the proofing engine would flag) sits in its own run,
#     matching how Word re-serializes the paragraph after it is edited/proofed.
#  2. Split the "Note: ... [TRAINER’S INITIALS]" run right after the opening
#     bracket, matching the grammar-check run boundary Word inserts there.
#  3. Split the " Slides" run into a leading space run and a "Slides" run, matching
#     the grammar-check run boundary Word inserts around that word.
#
# NOTE: Word marks these split points with <w:proofErr .../> elements that are
# purely cosmetic (no visible text, not exposed anywhere on the Range/Font object
# model) - they are stamped in automatically by Word's background spelling/grammar
# checker and are not reachable through COM automation. The run-splits that the
# proofErr elements straddle are reproduced here by nudging a Font property back
# to its current value, which is enough to stop Word from re-merging the runs.

$d = $word.ActiveDocument

# --- 1. "TDD Cookbook Development – CentOS 6.7 – 1.0.0" -> "... 1.1.0", split runs ---

# 1a. Bump the version number.
$d.Content.Find.Execute("1.0.0", $false, $false, $false, $false, $false, $true, 1, $false, "1.1.0", 2)

# 1b. Give "CentOS" its own run.
$rCentOS = $d.Content
$rCentOS.Find.Execute("CentOS", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$runCentOS = $rCentOS.Duplicate
$runCentOS.Font.Bold = $false
$runCentOS.Font.Bold = $true

# 1c. Give " 6.7 – " (between "CentOS" and the version number) its own run.
$rMid = $d.Content
$rMid.Find.Execute(" 6.7 – ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$runMid = $rMid.Duplicate
$runMid.Font.Bold = $false
$runMid.Font.Bold = $true

# --- 2. "Note: ...:  [TRAINER’S INITIALS]" - split right after "[" ---

$rBracket = $d.Content
$rBracket.Find.Execute(":  [", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$runBracket = $rBracket.Duplicate
$runBracket.Font.Italic = $false
$runBracket.Font.Italic = $true

# --- 3. " Slides" - split into " " and "Slides" ---

$rSlides = $d.Content
$rSlides.Find.Execute("Slides", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$runSlides = $rSlides.Duplicate
$runSlides.Font.Size = $runSlides.Font.Size + 2
$runSlides.Font.Size = $runSlides.Font.Size - 2
